# The first paragraph currently reads (across several runs, split up by a
# spell-check proofErr around the surname):
#   "I am Jesse Dachyshyn (who is very confused)! I am a Data Analytics student."
# Replace that whole sentence with a single new sentence, collapsing it into
# one run and leaving the trailing _GoBack bookmark (and everything else)
# untouched.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.Text = "I am Jesse Dachyshyn (who is very confused)! I am a Data Analytics student."
$find.Replacement.Text = "Welcome to the machine!"
$find.Execute(
    $find.Text,            # FindText
    $false,                # MatchCase
    $false,                # MatchWholeWord
    $false,                # MatchWildcards
    $false,                # MatchSoundsLike
    $false,                # MatchAllWordForms
    $true,                 # Forward
    1,                     # Wrap (wdFindContinue)
    $false,                # Format
    $find.Replacement.Text,# ReplaceWith
    2                      # Replace (wdReplaceAll)
)
